$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202 (shifts rows 202:451 down to 203:452)
$ws.Rows("202:202").Insert()

# Fill in the missing year 1800 and its value
$ws.Range("A202").Value = 1800
$ws.Range("B202").Value = 5.31101693276688

# Match the style of the other Year cells in column A (bordered/centered style)
$ws.Range("A202").Font.Bold = $true
$ws.Range("A202").HorizontalAlignment = $ws.Range("A201").HorizontalAlignment
$ws.Range("A202").VerticalAlignment = $ws.Range("A201").VerticalAlignment
$ws.Range("A202").Borders.LineStyle = $ws.Range("A201").Borders.LineStyle
